$d = $word.ActiveDocument

# Remove the last row of the (only) table, which contains the
# "#Colleen added / a chunk to test github and stuff" row.
$table = $d.Tables.Item(1)
$lastRow = $table.Rows.Last
$lastRow.Delete()

# Add a new paragraph after the table with a source-code style comment.
$p = $d.Paragraphs.Add($d.Range($table.Range.End, $table.Range.End))
$p.Range.Style = "SourceCode"
$commentText = "#Generating a table with percentages and 95% CIs by drug type for the years 1999 and 2010"
$p.Range.InsertAfter($commentText)

# Apply the CommentTok character style to the new run. Using
# Find/Replace (rather than Range.CharacterStyle, which mis-emits the
# formatting) ensures the proper run style (rStyle) is written.
$commentRange = $d.Range($p.Range.Start, $p.Range.End)
$commentRange.Find.ClearFormatting()
$commentRange.Find.Replacement.ClearFormatting()
$commentRange.Find.Replacement.Style = "CommentTok"
[void]$commentRange.Find.Execute($commentText, $false, $false, $false, $false, $false, $true, 1, $false, $commentText, 2)
